# Complete version of the presentation
# Including Statistical Modeling, Results, Conclusion, and Future Work

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Fix a typo on the title slide (slide 1): "Feng ,Tianrui" -> "Feng, Tianrui"
#    Do this BEFORE duplicating the slide so the new closing slide
#    inherits the corrected text too.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$authorsShape = $s1.Shapes.Item(3)
$authorsRange = $authorsShape.TextFrame.TextRange
$firstRun = $authorsRange.Characters(1, 48)
$firstRun.Text = "Jia Niu, Xiaolin Feng, Tianrui Xie, Zhilin Liu, "

# ---------------------------------------------------------------------
# 2. Add a new closing "Thanks for your attention!" slide at the end
#    (position 12), built from a duplicate of slide 1 (title slide),
#    matching its layout/background/picture/author block.
# ---------------------------------------------------------------------
$dup = $s1.Duplicate()
$newSlide = $dup.Item(1)
$newSlide.MoveTo($p.Slides.Count)

# -- Resize/reposition + retext the "Title 2" placeholder --------------
$titleShape = $newSlide.Shapes.Item(2)
$titleRange = $titleShape.TextFrame.TextRange
$oldFirstRun = $titleRange.Characters(1, 9)
$oldFirstRun.Text = "Thanks for your attention!"
$newLen = "Thanks for your attention!".Length
$tailLen = $titleRange.Length - $newLen
if ($tailLen -gt 0) {
    $tail = $titleRange.Characters($newLen + 1, $tailLen)
    $tail.Delete()
}

# Reapply formatting on the remaining (now sole) run: bigger size, no shadow
$runRange = $titleRange.Characters(1, $newLen)
$runRange.Font.Size = 40
$runRange.Font.Shadow = $false

# Move/resize the title placeholder to its new position
$titleShape.Left = 762000 / 914400 * 72
$titleShape.Top = 1824395 / 914400 * 72
$titleShape.Width = 7772400 / 914400 * 72
$titleShape.Height = 1107996 / 914400 * 72

Write-Output "done"
